# Add 100v90 and 100v95 RAREsim v2.1.1 results (rows 54-59)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$note90 = "1. Separately-RAREsim v2.1.1`n--functional_bins 100%`n--synonymous_bins 100%`n2. RAREsim v2.1.1: Convert 100% pruned hap file to .sm file using convert.py`n3. Separately-RAREsim v2.1.1`n--functional_bins 90% 6 MAC BINS`n--synonymous_bins 90% 6 MAC BINS`n4. R: add pruned variants back in as rows of 0 and extract datasets"

$note95 = "1. Separately-RAREsim v2.1.1`n--functional_bins 100%`n--synonymous_bins 100%`n2. RAREsim v2.1.1: Convert 100% pruned hap file to .sm file using convert.py`n3. Separately-RAREsim v2.1.1`n--functional_bins 95% 6 MAC BINS`n--synonymous_bins 95% 6 MAC BINS`n4. R: add pruned variants back in as rows of 0 and extract datasets"

$cols = @("J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH")

# Row 54: internal 1, external 1  (100v90 pipeline)
$row54 = @(0.03,0.05,0.03,0.04,0.05,0.03,0.01,0.06,0.04,0.04,0.03,0.01,0.01,0.02,0.02,0.04,0.06,0.07,0.04,0,0,0.01,0.05,0.05,0.03)
# Row 55: internal 1, external 0.9  (100v90 pipeline)
$row55 = @(0.16,0.05,0.16,0.04,0.05,0.03,0.01,0.68,0.43,0.04,0.03,0.49,0.34,0.02,0.56,0.25,0.06,0.59,0.39,0,0.45,0.26,0.05,0.46,0.27)
# Row 56: internal 0.9, external 0.9  (100v90 pipeline)
$row56 = @(0.04,0.05,0.04,0.06,0.05,0.03,0.03,0.05,0.03,0.03,0.02,0.02,0.02,0.01,0.05,0.03,0.02,0.05,0.03,0,0,0.02,0.05,0.06,0.05)
# Row 57: internal 1, external 1  (100v95 pipeline)
$row57 = @(0.06,0.07,0.06,0.07,0.07,0,0.01,0,0.02,0.01,0.01,0.04,0.03,0.01,0,0.01,0.04,0.05,0.03,0.01,0,0.01,0.08,0.05,0.03)
# Row 58: internal 1, external 0.95  (100v95 pipeline)
$row58 = @(0.14,0.07,0.14,0.07,0.07,0.01,0.01,0.45,0.25,0.02,0.01,0.29,0.17,0.01,0.4,0.23,0.04,0.36,0.22,0.01,0.29,0.14,0.08,0.23,0.12)
# Row 59: internal 0.95, external 0.95  (100v95 pipeline)
$row59 = @(0.05,0.08,0.05,0.08,0.08,0,0.01,0.03,0.03,0.02,0.02,0.03,0.02,0.01,0,0,0.06,0.06,0.06,0.01,0.01,0.01,0.04,0.06,0.03)

$rowsData = @($row54, $row55, $row56, $row57, $row58, $row59)

$rowNums = @(54, 55, 56, 57, 58, 59)
$dVals   = @(1,    1,    0.9,  1,    1,    0.95)
$eVals   = @(1,    0.9,  0.9,  1,    0.95, 0.95)
$fVals   = @("100v90", "100v90", "100v90", "100v95", "100v95", "100v95")
$iVals   = @($note90, $note90, $note90, $note95, $note95, $note95)

# 1) Bring formatting for each new row in line with the existing table rows
#    by copying the formats from row 53 (same pattern used for previous entries).
for ($idx = 0; $idx -lt $rowNums.Length; $idx++) {
    $r = $rowNums[$idx]
    $ws.Range("A53:AH53").Copy()
    $ws.Range("A" + $r + ":AH" + $r).PasteSpecial(-4122) # xlPasteFormats
}
$ws.Application.CutCopyMode = $false

# 2) Fill in the values for each new row
for ($idx = 0; $idx -lt $rowNums.Length; $idx++) {
    $r = $rowNums[$idx]

    $ws.Range("A$r").Value = 45219
    $ws.Range("B$r").Value = "100% NFE"
    $ws.Range("C$r").Value = 0.001
    $ws.Range("D$r").Value = $dVals[$idx]
    $ws.Range("E$r").Value = $eVals[$idx]
    $ws.Range("F$r").Value = $fVals[$idx]
    $ws.Range("G$r").Value = "RAREsim v2.1.1"
    $ws.Range("H$r").Value = "Hapgen Haplotype"
    $ws.Range("I$r").Value = $iVals[$idx]

    $vals = $rowsData[$idx]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + "$r").Value = $vals[$c]
    }

    $ws.Rows.Item($r).RowHeight = 158.4
}

# 3) Update the selection to match the post-edit state
$ws.Range("J61").Select()
